$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scraped data for rows 2-9 (row 9 is a brand-new row; rows 2-8 are overwritten)
$data = @(
    @("1331155", "https://aiesec.org/opportunity/global-talent/1331155", "Co-Manufacturing Trainee", "Panamá, Provincia de Panamá, Panamá", "No", "0 applicants", "6 - 18 Months", "NESTLE"),
    @("1331153", "https://aiesec.org/opportunity/global-talent/1331153", "[Impact Brazil] - GTM Engineer Intern", "São Paulo, SP, Brasil", "No", "2 applicants", "3 - 6 Months", "Ecomiles"),
    @("1331110", "https://aiesec.org/opportunity/global-talent/1331110", "Content Creator", "Sousse, Tunisie", "No", "1 applicant", "9 - 12 Weeks", "Next Round"),
    @("1331109", "https://aiesec.org/opportunity/global-talent/1331109", "Web Developer", "Sousse, Tunisie", "No", "1 applicant", "9 - 12 Weeks", "Next Round"),
    @("1331101", "https://aiesec.org/opportunity/global-talent/1331101", "Mobile Developer", "Sousse, Tunisie", "No", "0 applicants", "9 - 12 Weeks", "Business 360"),
    @("1331055", "https://aiesec.org/opportunity/global-talent/1331055", "Repair technician", "Hammam Sousse, Tunisie", "No", "0 applicants", "3 - 6 Months", "MOBYSTORE"),
    @("1328206", "https://aiesec.org/opportunity/global-talent/1328206", "Power BI Specialist", "Frankfurt am Main, Deutschland", "No", "193 applicants", "3 - 6 Months", "Greyfood GmbH"),
    @("1321497", "https://aiesec.org/opportunity/global-talent/1321497", "Sales  Specialist", "Kartepe, Kocaeli, Türkiye", "No", "72 applicants", "6 - 18 Months", "Dessa Teknoloji Sanayi Ticaret Limited Şirketi")
)

# Column A holds numeric-looking opportunity IDs that must stay text, like the rest
# of the sheet (matches the source data's text formatting instead of auto-converting
# to numbers).
$ws.Range("A2:A9").NumberFormat = "@"

# E2 previously carried a highlighted ("Yes"/premium) style; the new value is a plain
# "No" so the old yellow-fill formatting must be dropped back to Normal first.
$ws.Range("E2").Style = "Normal"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# Column width updates (C: 66->40, D: 49->38, H: 45->49). ColumnWidth needs a small
# offset subtracted to land exactly on the target stored width.
$ws.Columns.Item(3).ColumnWidth = 40 - 5/6
$ws.Columns.Item(4).ColumnWidth = 38 - 5/6
$ws.Columns.Item(8).ColumnWidth = 49 - 5/6
